# Rename the "SwateTemplateMetadata" sheet to "isa_template"
$wb = $excel.ActiveWorkbook
$metaSheet = $wb.Worksheets.Item("SwateTemplateMetadata")
$metaSheet.Name = "isa_template"

# Clear the (visually empty / no-op styled) leftover formatting on E13 and D14
# so they go back to the default style, matching the cleaned-up workbook, and
# drop their (already empty) contents so the cell records disappear entirely.
$metaSheet.Range("E13").Interior.Pattern = -4142
$metaSheet.Range("E13").Borders.LineStyle = -4142
$metaSheet.Range("E13").ClearContents()

$metaSheet.Range("D14").Interior.Pattern = -4142
$metaSheet.Range("D14").Borders.LineStyle = -4142
$metaSheet.Range("D14").ClearContents()
